$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 25
$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 25

$ws.Range("A6").Select()
